# "Updated: po 22. 02. 2021"
# Revises the AgTests (F) / AgPosit (G) figures for a span of previously
# reported days, and appends three new daily rows (2021-02-19 .. 02-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised AgTests / AgPosit values for existing rows ---
$ws.Cells.Item(281, 6).Value = 45603
$ws.Cells.Item(282, 6).Value = 47262
$ws.Cells.Item(285, 6).Value = 41929
$ws.Cells.Item(285, 7).Value = 3431
$ws.Cells.Item(286, 6).Value = 55258
$ws.Cells.Item(292, 6).Value = 82063
$ws.Cells.Item(292, 7).Value = 7273
$ws.Cells.Item(293, 6).Value = 82104
$ws.Cells.Item(293, 7).Value = 5762
$ws.Cells.Item(294, 6).Value = 93091
$ws.Cells.Item(294, 7).Value = 4891
$ws.Cells.Item(299, 6).Value = 65133
$ws.Cells.Item(299, 7).Value = 6884
$ws.Cells.Item(300, 6).Value = 71910
$ws.Cells.Item(300, 7).Value = 7057
$ws.Cells.Item(301, 6).Value = 71359
$ws.Cells.Item(301, 7).Value = 5640
$ws.Cells.Item(302, 6).Value = 77386
$ws.Cells.Item(302, 7).Value = 5713
$ws.Cells.Item(306, 6).Value = 71707
$ws.Cells.Item(306, 7).Value = 7218
$ws.Cells.Item(307, 6).Value = 75622
$ws.Cells.Item(307, 7).Value = 6443
$ws.Cells.Item(309, 6).Value = 75146
$ws.Cells.Item(309, 7).Value = 5318
$ws.Cells.Item(310, 6).Value = 75909
$ws.Cells.Item(310, 7).Value = 3956
$ws.Cells.Item(313, 6).Value = 72288
$ws.Cells.Item(313, 7).Value = 3260
$ws.Cells.Item(314, 6).Value = 63885
$ws.Cells.Item(314, 7).Value = 3170
$ws.Cells.Item(315, 6).Value = 56207
$ws.Cells.Item(315, 7).Value = 2650
$ws.Cells.Item(316, 6).Value = 49685
$ws.Cells.Item(316, 7).Value = 2246
$ws.Cells.Item(317, 6).Value = 62241
$ws.Cells.Item(317, 7).Value = 2147
$ws.Cells.Item(320, 6).Value = 69383
$ws.Cells.Item(320, 7).Value = 3154
$ws.Cells.Item(321, 6).Value = 93524
$ws.Cells.Item(321, 7).Value = 2843
$ws.Cells.Item(322, 6).Value = 107333
$ws.Cells.Item(322, 7).Value = 2302
$ws.Cells.Item(323, 6).Value = 214186
$ws.Cells.Item(323, 7).Value = 3175
$ws.Cells.Item(324, 6).Value = 235470
$ws.Cells.Item(324, 7).Value = 2718
$ws.Cells.Item(325, 6).Value = 763192
$ws.Cells.Item(325, 7).Value = 6498
$ws.Cells.Item(326, 6).Value = 433560
$ws.Cells.Item(326, 7).Value = 3848
$ws.Cells.Item(327, 6).Value = 237967
$ws.Cells.Item(327, 7).Value = 2900
$ws.Cells.Item(328, 6).Value = 180664
$ws.Cells.Item(328, 7).Value = 2649
$ws.Cells.Item(329, 6).Value = 89585
$ws.Cells.Item(329, 7).Value = 1818
$ws.Cells.Item(330, 6).Value = 71005
$ws.Cells.Item(330, 7).Value = 2004
$ws.Cells.Item(331, 6).Value = 151381
$ws.Cells.Item(331, 7).Value = 2639
$ws.Cells.Item(332, 6).Value = 436828
$ws.Cells.Item(332, 7).Value = 4306
$ws.Cells.Item(333, 6).Value = 265797
$ws.Cells.Item(333, 7).Value = 2859
$ws.Cells.Item(334, 6).Value = 202541
$ws.Cells.Item(334, 7).Value = 3372
$ws.Cells.Item(335, 6).Value = 129555
$ws.Cells.Item(335, 7).Value = 2906
$ws.Cells.Item(336, 6).Value = 101268
$ws.Cells.Item(336, 7).Value = 3212
$ws.Cells.Item(337, 6).Value = 102801
$ws.Cells.Item(337, 7).Value = 2923
$ws.Cells.Item(338, 6).Value = 220000
$ws.Cells.Item(338, 7).Value = 3079
$ws.Cells.Item(339, 6).Value = 644627
$ws.Cells.Item(339, 7).Value = 5503
$ws.Cells.Item(340, 6).Value = 381353
$ws.Cells.Item(340, 7).Value = 3280
$ws.Cells.Item(341, 6).Value = 294810
$ws.Cells.Item(341, 7).Value = 3653
$ws.Cells.Item(342, 6).Value = 173810
$ws.Cells.Item(342, 7).Value = 2948
$ws.Cells.Item(343, 6).Value = 127567
$ws.Cells.Item(343, 7).Value = 2833
$ws.Cells.Item(344, 6).Value = 130880
$ws.Cells.Item(344, 7).Value = 2416
$ws.Cells.Item(345, 6).Value = 279323
$ws.Cells.Item(345, 7).Value = 3175
$ws.Cells.Item(346, 6).Value = 644069
$ws.Cells.Item(346, 7).Value = 4575
$ws.Cells.Item(347, 6).Value = 328502
$ws.Cells.Item(347, 7).Value = 2760
$ws.Cells.Item(348, 6).Value = 224889
$ws.Cells.Item(348, 7).Value = 3078
$ws.Cells.Item(349, 6).Value = 161668
$ws.Cells.Item(349, 7).Value = 2720
$ws.Cells.Item(350, 6).Value = 120708
$ws.Cells.Item(350, 7).Value = 2610
$ws.Cells.Item(351, 6).Value = 138817
$ws.Cells.Item(351, 7).Value = 2610

# --- New daily rows appended at the bottom (352-354) ---
$ws.Cells.Item(352, 1).Value = 44246
$ws.Cells.Item(352, 2).Value = 290457
$ws.Cells.Item(352, 3).Value = 14300
$ws.Cells.Item(352, 4).Value = 2705
$ws.Cells.Item(352, 5).Value = 6424
$ws.Cells.Item(352, 6).Value = 270721
$ws.Cells.Item(352, 7).Value = 3567

$ws.Cells.Item(353, 1).Value = 44247
$ws.Cells.Item(353, 2).Value = 292143
$ws.Cells.Item(353, 3).Value = 8277
$ws.Cells.Item(353, 4).Value = 1686
$ws.Cells.Item(353, 5).Value = 6505
$ws.Cells.Item(353, 6).Value = 566613
$ws.Cells.Item(353, 7).Value = 4392

$ws.Cells.Item(354, 1).Value = 44248
$ws.Cells.Item(354, 2).Value = 292792
$ws.Cells.Item(354, 3).Value = 3059
$ws.Cells.Item(354, 4).Value = 649
$ws.Cells.Item(354, 5).Value = 6577
$ws.Cells.Item(354, 6).Value = 228403
$ws.Cells.Item(354, 7).Value = 2384

# Ensure the new date cells pick up the same date style (yyyy-mm-dd) as the
# rest of column A, in case style inheritance from the column default
# doesn't apply for some reason.
$ws.Range("A352:A354").NumberFormat = "yyyy-mm-dd"
